# Auto-generated Excel COM-interop script
# Updates H:N market-data columns across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW)
# to reflect refreshed Leve profit calculations from the scheduled data-fetch runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 90.40000000000001
$ws.Range("J2").Value = 92
$ws.Range("L2").Value = 92
$ws.Range("N2").Value = -318
$ws.Range("H5").Value = 110.25
$ws.Range("I5").Value = 72.5
$ws.Range("J5").Value = 148
$ws.Range("K5").Value = 72.5
$ws.Range("L5").Value = 148
$ws.Range("M5").Value = 42.5
$ws.Range("N5").Value = -378
$ws.Range("H6").Value = 20223500
$ws.Range("I6").Value = 10670000
$ws.Range("J6").Value = 25000250
$ws.Range("K6").Value = 32010000
$ws.Range("L6").Value = 75000750
$ws.Range("M6").Value = -32009888
$ws.Range("N6").Value = -75000974
$ws.Range("H9").Value = 254.57143
$ws.Range("I9").Value = 100
$ws.Range("J9").Value = 370.5
$ws.Range("K9").Value = 100
$ws.Range("L9").Value = 370.5
$ws.Range("M9").Value = 69
$ws.Range("N9").Value = -708.5
$ws.Range("H10").Value = 4995
$ws.Range("J10").Value = 4995
$ws.Range("L10").Value = 4995
$ws.Range("N10").Value = -5581
$ws.Range("H12").Value = 424
$ws.Range("I12").Value = 380
$ws.Range("J12").Value = 453.33334
$ws.Range("K12").Value = 380
$ws.Range("L12").Value = 453.33334
$ws.Range("M12").Value = -210
$ws.Range("N12").Value = -793.33334
$ws.Range("H16").Value = 2500
$ws.Range("J16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2460
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = $null
$ws.Range("N18").Value = $null
$ws.Range("H19").Value = 221.1
$ws.Range("I19").Value = 167.4
$ws.Range("J19").Value = 274.8
$ws.Range("K19").Value = 167.4
$ws.Range("L19").Value = 274.8
$ws.Range("M19").Value = 7.599999999999994
$ws.Range("N19").Value = -624.8
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = $null
$ws.Range("H21").Value = 3999
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = $null
$ws.Range("H23").Value = 3999
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = $null
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = $null
$ws.Range("N29").Value = $null
$ws.Range("H31").Value = 928.5714
$ws.Range("I31").Value = 700
$ws.Range("J31").Value = 1500
$ws.Range("K31").Value = 2100
$ws.Range("L31").Value = 4500
$ws.Range("M31").Value = -1870
$ws.Range("N31").Value = -4960
$ws.Range("H32").Value = 1028.7333
$ws.Range("I32").Value = 725
$ws.Range("J32").Value = 1075.4615
$ws.Range("K32").Value = 725
$ws.Range("L32").Value = 1075.4615
$ws.Range("M32").Value = -399
$ws.Range("N32").Value = -1727.4615
$ws.Range("H33").Value = 6428.6875
$ws.Range("J33").Value = 17033
$ws.Range("L33").Value = 17033
$ws.Range("N33").Value = -17491
$ws.Range("H34").Value = 2660
$ws.Range("I34").Value = 1546.6666
$ws.Range("K34").Value = 1546.6666
$ws.Range("M34").Value = -1343.6666
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = $null
$ws.Range("H36").Value = 2660
$ws.Range("I36").Value = 1546.6666
$ws.Range("K36").Value = 1546.6666
$ws.Range("M36").Value = -831.6666
$ws.Range("H38").Value = 121132.68
$ws.Range("I38").Value = 166797.61
$ws.Range("J38").Value = 3708.5715
$ws.Range("K38").Value = 500392.83
$ws.Range("L38").Value = 11125.7145
$ws.Range("M38").Value = -500020.83
$ws.Range("N38").Value = -11869.7145
$ws.Range("H39").Value = 126.72222
$ws.Range("I39").Value = 62.22222
$ws.Range("J39").Value = 191.22223
$ws.Range("K39").Value = 186.66666
$ws.Range("L39").Value = 573.66669
$ws.Range("M39").Value = 109.33334
$ws.Range("N39").Value = -1165.66669
$ws.Range("H40").Value = 1005.9
$ws.Range("I40").Value = 870.2
$ws.Range("J40").Value = 1141.6
$ws.Range("K40").Value = 870.2
$ws.Range("L40").Value = 1141.6
$ws.Range("M40").Value = -695.2
$ws.Range("N40").Value = -1491.6
$ws.Range("H41").Value = 1819.1428
$ws.Range("I41").Value = 2800
$ws.Range("J41").Value = 511.33334
$ws.Range("K41").Value = 2800
$ws.Range("L41").Value = 511.33334
$ws.Range("M41").Value = -2360
$ws.Range("N41").Value = -1391.33334
$ws.Range("H42").Value = 109.90909
$ws.Range("I42").Value = 10
$ws.Range("J42").Value = 193.16667
$ws.Range("K42").Value = 30
$ws.Range("L42").Value = 579.50001
$ws.Range("M42").Value = 200
$ws.Range("N42").Value = -1039.50001
$ws.Range("H43").Value = 925.0454999999999
$ws.Range("I43").Value = 1085.5
$ws.Range("J43").Value = 833.3570999999999
$ws.Range("K43").Value = 1085.5
$ws.Range("L43").Value = 833.3570999999999
$ws.Range("M43").Value = -1016.5
$ws.Range("N43").Value = -971.3570999999999
$ws.Range("H44").Value = 95000
$ws.Range("J44").Value = 95000
$ws.Range("L44").Value = 95000
$ws.Range("N44").Value = -95924
$ws.Range("H46").Value = 1115.3636
$ws.Range("I46").Value = 1270
$ws.Range("J46").Value = 1027
$ws.Range("K46").Value = 3810
$ws.Range("L46").Value = 3081
$ws.Range("M46").Value = -3691
$ws.Range("N46").Value = -3319
$ws.Range("H47").Value = 7600
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 7600
$ws.Range("K47").Value = 0
$ws.Range("M47").Value = $null
$ws.Range("N47").Value = -9544
$ws.Range("H48").Value = 2250
$ws.Range("J48").Value = 2250
$ws.Range("L48").Value = 6750
$ws.Range("N48").Value = -7334
$ws.Range("H51").Value = 3708
$ws.Range("J51").Value = 3724
$ws.Range("L51").Value = 3724
$ws.Range("N51").Value = -4692
$ws.Range("H54").Value = 4284
$ws.Range("I54").Value = 2978.6667
$ws.Range("J54").Value = 8200
$ws.Range("K54").Value = 2978.6667
$ws.Range("L54").Value = 8200
$ws.Range("M54").Value = -2492.6667
$ws.Range("N54").Value = -9172
$ws.Range("H56").Value = 2250
$ws.Range("J56").Value = 2250
$ws.Range("L56").Value = 6750
$ws.Range("N56").Value = -7818
$ws.Range("H58").Value = 2622.5
$ws.Range("I58").Value = 1500
$ws.Range("J58").Value = 2996.6667
$ws.Range("K58").Value = 4500
$ws.Range("L58").Value = 8990.000100000001
$ws.Range("M58").Value = -4350
$ws.Range("N58").Value = -9290.000100000001
$ws.Range("H59").Value = 1361.6
$ws.Range("J59").Value = 1677
$ws.Range("L59").Value = 5031
$ws.Range("N59").Value = -6145
$ws.Range("H60").Value = 1115.3636
$ws.Range("I60").Value = 1270
$ws.Range("J60").Value = 1027
$ws.Range("K60").Value = 3810
$ws.Range("L60").Value = 3081
$ws.Range("M60").Value = -3326
$ws.Range("N60").Value = -4049
$ws.Range("H61").Value = 489
$ws.Range("I61").Value = 431.4
$ws.Range("J61").Value = 777
$ws.Range("K61").Value = 1294.2
$ws.Range("L61").Value = 2331
$ws.Range("M61").Value = -1122.2
$ws.Range("N61").Value = -2675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2883.9092
$ws.Range("I74").Value = 2957.4614
$ws.Range("J74").Value = 2777.6667
$ws.Range("K74").Value = 2957.4614
$ws.Range("L74").Value = 2777.6667
$ws.Range("M74").Value = -2083.4614
$ws.Range("N74").Value = -4525.6667
$ws.Range("H77").Value = 2883.9092
$ws.Range("I77").Value = 2957.4614
$ws.Range("J77").Value = 2777.6667
$ws.Range("K77").Value = 14787.307
$ws.Range("L77").Value = 13888.3335
$ws.Range("M77").Value = -10419.307
$ws.Range("N77").Value = -22624.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 44113.332
$ws.Range("J122").Value = 44113.332
$ws.Range("L122").Value = 44113.332
$ws.Range("N122").Value = -53913.332
$ws.Range("H123").Value = 36527
$ws.Range("J123").Value = 36527
$ws.Range("L123").Value = 36527
$ws.Range("N123").Value = -46327
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820
$ws.Range("H125").Value = 33000
$ws.Range("J125").Value = 33000
$ws.Range("L125").Value = 33000
$ws.Range("N125").Value = -42840
$ws.Range("H126").Value = 44389.5
$ws.Range("J126").Value = 44389.5
$ws.Range("L126").Value = 44389.5
$ws.Range("N126").Value = -54269.5
$ws.Range("H127").Value = 69770
$ws.Range("J127").Value = 69770
$ws.Range("L127").Value = 69770
$ws.Range("N127").Value = -79690
$ws.Range("H129").Value = 48872
$ws.Range("J129").Value = 48872
$ws.Range("L129").Value = 48872
$ws.Range("N129").Value = -58872
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("N130").Value = $null
$ws.Range("H131").Value = 54780
$ws.Range("J131").Value = 54780
$ws.Range("L131").Value = 54780
$ws.Range("N131").Value = -64860
$ws.Range("H140").Value = 65650
$ws.Range("J140").Value = 65650
$ws.Range("L140").Value = 65650
$ws.Range("N140").Value = -76010

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 15296.667
$ws.Range("J63").Value = 15296.667
$ws.Range("L63").Value = 15296.667
$ws.Range("N63").Value = -16668.667
$ws.Range("H66").Value = 15296.667
$ws.Range("J66").Value = 15296.667
$ws.Range("L66").Value = 45890.001
$ws.Range("N66").Value = -52754.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 324.0909
$ws.Range("I33").Value = 101.9375
$ws.Range("J33").Value = 916.5
$ws.Range("K33").Value = 611.625
$ws.Range("L33").Value = 5499
$ws.Range("M33").Value = -328.625
$ws.Range("N33").Value = -6065

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2390.5
$ws.Range("I122").Value = 2391.7144
$ws.Range("J122").Value = 2387.6667
$ws.Range("K122").Value = 7175.1432
$ws.Range("L122").Value = 7163.000100000001
$ws.Range("M122").Value = -4725.1432
$ws.Range("N122").Value = -12063.0001
$ws.Range("H132").Value = 1926.7179
$ws.Range("I132").Value = 1471
$ws.Range("J132").Value = 3445.7778
$ws.Range("K132").Value = 4413
$ws.Range("L132").Value = 10337.3334
$ws.Range("M132").Value = -1883
$ws.Range("N132").Value = -15397.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 780.4400000000001
$ws.Range("I55").Value = 1106.4166
$ws.Range("J55").Value = 479.53845
$ws.Range("K55").Value = 1106.4166
$ws.Range("L55").Value = 479.53845
$ws.Range("M55").Value = -933.4166
$ws.Range("N55").Value = -825.53845

Write-Host "Applied scheduled market-data refresh to ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets."
